# Insert a new weekly price-report row for "Apio" (Vega Modelo de Temuco)
# above the existing row 581, shifting all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 581 (pushes old rows 581..612 to 582..613).
$ws.Rows.Item(581).Insert()

# Populate the newly inserted row 581 with the new weekly data point.
$ws.Cells.Item(581, 1).Value = 10
$ws.Cells.Item(581, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(581, 3).Value = "La Araucanía"
$ws.Cells.Item(581, 4).Value = 45267
$ws.Cells.Item(581, 5).Value = 9
$ws.Cells.Item(581, 6).Value = 100112017
$ws.Cells.Item(581, 7).Value = "Apio"
$ws.Cells.Item(581, 8).Value = "Americana (o)"
$ws.Cells.Item(581, 9).Value = "Primera"
$ws.Cells.Item(581, 10).Value = 600
$ws.Cells.Item(581, 11).Value = 11000
$ws.Cells.Item(581, 12).Value = 11000
$ws.Cells.Item(581, 13).Value = 11000
$ws.Cells.Item(581, 14).Value = "$/caja 8 unidades"
$ws.Cells.Item(581, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(581, 16).Value = 11000
$ws.Cells.Item(581, 17).Value = 1
$ws.Cells.Item(581, 18).Value = "Hortaliza"
